# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps for the a20803dc-... report row,
# as part of regenerating the handback status report.

$wb = $excel.ActiveWorkbook

# Overview sheet: G3 holds the "Latest HO Xliff Generate Date" for the
# a20803dc-ec67-44aa-ae94-89315acf7b70.md row.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-30 20:51:22"

# zh-cn sheet: H3 is "Correspond Handoff Datetime", K3 is
# "Correspond Handback DateTime" for the same row.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-30 20:51:18"
$wsZhCn.Range("K3").Value = "2016-08-30 20:51:35"

# de-de sheet: K3 is "Correspond Handback DateTime" for the same row.
# (de-de's H3 shares the same string as Overview!G3 and is updated
# automatically through the shared string table.)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K3").Value = "2016-08-30 20:51:42"
